$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Value = "X"
$ws.Range("C18").Value = "X"
$ws.Range("D18").Value = "X"
$ws.Range("E18").Value = "X"
$ws.Range("F18").Value = -20
$ws.Range("G18").Value = 52
$ws.Range("I18").Value = 256
$ws.Range("J18").Value = -220
$ws.Range("K18").Value = "X"

$ws.Range("M21").Select()
